# Add I0 and IF columns (I and J) to the sheet, mirroring the style of
# the existing header row and filling in the computed values for each
# data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers: I1 = "I0", J1 = "IF" with same style as existing header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$values = @{
    2  = @(5, 6)
    3  = @(5, 6)
    4  = @(6, 7)
    5  = @(6, 6)
    6  = @(7, 7)
    7  = @(6, 7)
    8  = @(4, 6)
    9  = @(9, 9)
    10 = @(6, 6)
    11 = @(3, 4)
    12 = @(6, 6)
    13 = @(1, 3)
    14 = @(5, 6)
    15 = @(6, 7)
    16 = @(6, 7)
    17 = @(7, 7)
    18 = @(6, 7)
    19 = @(5, 6)
    20 = @(6, 6)
    21 = @(4, 4)
    22 = @(8, 8)
    23 = @(8, 8)
    24 = @(12, 12)
    25 = @(4, 5)
    26 = @(6, 6)
    27 = @(6, 7)
    28 = @(6, 7)
    29 = @(7, 7)
    30 = @(6, 7)
    31 = @(5, 5)
    32 = @(8, 8)
    33 = @(7, 7)
    34 = @(8, 9)
    35 = @(7, 7)
    36 = @(8, 8)
    37 = @(9, 9)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
